$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1806.1818
$ws.Range("I19").Value = 699.8
$ws.Range("K19").Value = 699.8
$ws.Range("M19").Value = -524.8
$ws.Range("H64").Value = 3260.8696
$ws.Range("J64").Value = 6000
$ws.Range("L64").Value = 6000
$ws.Range("N64").Value = -6496
$ws.Range("H67").Value = 3260.8696
$ws.Range("J67").Value = 6000
$ws.Range("L67").Value = 6000
$ws.Range("N67").Value = -7716
$ws.Range("H74").Value = 5292.3076
$ws.Range("I74").Value = 4987
$ws.Range("K74").Value = 4987
$ws.Range("M74").Value = -4051
$ws.Range("H77").Value = 5292.3076
$ws.Range("I77").Value = 4987
$ws.Range("K77").Value = 24935
$ws.Range("M77").Value = -20255
$ws.Range("H98").Value = 1966
$ws.Range("I98").Value = 1400
$ws.Range("J98").Value = 2249
$ws.Range("K98").Value = 1400
$ws.Range("L98").Value = 2249
$ws.Range("M98").Value = 98
$ws.Range("N98").Value = -5245
$ws.Range("H122").Value = 1966
$ws.Range("I122").Value = 1400
$ws.Range("J122").Value = 2249
$ws.Range("K122").Value = 4200
$ws.Range("L122").Value = 6747
$ws.Range("M122").Value = -1750
$ws.Range("N122").Value = -11647
$ws.Range("H133").Value = 74047.82000000001
$ws.Range("J133").Value = 74047.82000000001
$ws.Range("L133").Value = 74047.82000000001
$ws.Range("N133").Value = -84167.82000000001
$ws.Range("H134").Value = 99999
$ws.Range("J134").Value = 99999
$ws.Range("L134").Value = 99999
$ws.Range("N134").Value = -110139
$ws.Range("H135").Value = 1686.6666
$ws.Range("I135").Value = 1486.0667
$ws.Range("K135").Value = 13374.6003
$ws.Range("M135").Value = -10839.6003
$ws.Range("H136").Value = 99995
$ws.Range("J136").Value = 99995
$ws.Range("L136").Value = 99995
$ws.Range("N136").Value = -110195
$ws.Range("H137").Value = 404838.9
$ws.Range("J137").Value = 727049.3
$ws.Range("L137").Value = 2181147.9
$ws.Range("N137").Value = -2186247.9
$ws.Range("H139").Value = 98406
$ws.Range("J139").Value = 98406
$ws.Range("L139").Value = 98406
$ws.Range("N139").Value = -108686
$ws.Range("H140").Value = 80776.664
$ws.Range("J140").Value = 80776.664
$ws.Range("L140").Value = 80776.664
$ws.Range("N140").Value = -91136.664
$ws.Range("H141").Value = 4924.773
$ws.Range("I141").Value = 3962.1177
$ws.Range("K141").Value = 11886.3531
$ws.Range("M141").Value = -6706.3531

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 697.5714
$ws.Range("I61").Value = 645.3333
$ws.Range("J61").Value = 1011
$ws.Range("K61").Value = 645.3333
$ws.Range("L61").Value = 1011
$ws.Range("M61").Value = -433.3333
$ws.Range("N61").Value = -1435
$ws.Range("H74").Value = 2107.5
$ws.Range("I74").Value = 1564.8334
$ws.Range("K74").Value = 1564.8334
$ws.Range("M74").Value = -690.8334
$ws.Range("H77").Value = 2107.5
$ws.Range("I77").Value = 1564.8334
$ws.Range("K77").Value = 7824.166999999999
$ws.Range("M77").Value = -3456.166999999999
$ws.Range("H122").Value = 3332.6667
$ws.Range("I122").Value = 3332.6667
$ws.Range("K122").Value = 9998.000100000001
$ws.Range("M122").Value = -7548.000100000001
$ws.Range("H130").Value = 39900
$ws.Range("J130").Value = 39900
$ws.Range("L130").Value = 39900
$ws.Range("N130").Value = -49940
$ws.Range("H132").Value = 1532.9412
$ws.Range("I132").Value = 1249.75
$ws.Range("J132").Value = 2212.6
$ws.Range("K132").Value = 3749.25
$ws.Range("L132").Value = 6637.799999999999
$ws.Range("M132").Value = -1219.25
$ws.Range("N132").Value = -11697.8
$ws.Range("H134").Value = 124429
$ws.Range("J134").Value = 124429
$ws.Range("L134").Value = 124429
$ws.Range("N134").Value = -134569
$ws.Range("H136").Value = 697.5714
$ws.Range("I136").Value = 645.3333
$ws.Range("J136").Value = 1011
$ws.Range("K136").Value = 1935.9999
$ws.Range("L136").Value = 3033
$ws.Range("M136").Value = 614.0001
$ws.Range("N136").Value = -8133
$ws.Range("H138").Value = 150000
$ws.Range("J138").Value = 150000
$ws.Range("L138").Value = 150000
$ws.Range("N138").Value = -160280

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 19335.334
$ws.Range("I82").Value = 9003.5
$ws.Range("J82").Value = 39999
$ws.Range("K82").Value = 9003.5
$ws.Range("L82").Value = 39999
$ws.Range("M82").Value = -8620.5
$ws.Range("N82").Value = -40765
$ws.Range("H85").Value = 19335.334
$ws.Range("I85").Value = 9003.5
$ws.Range("J85").Value = 39999
$ws.Range("K85").Value = 9003.5
$ws.Range("L85").Value = 39999
$ws.Range("M85").Value = -7677.5
$ws.Range("N85").Value = -42651
$ws.Range("H107").Value = 2472.4443
$ws.Range("I107").Value = 1949.5
$ws.Range("K107").Value = 1949.5
$ws.Range("M107").Value = -29.5
$ws.Range("H134").Value = 3040
$ws.Range("I134").Value = 2415.8975
$ws.Range("K134").Value = 7247.6925
$ws.Range("M134").Value = -4712.6925
$ws.Range("H138").Value = 99769.42999999999
$ws.Range("J138").Value = 99769.42999999999
$ws.Range("L138").Value = 99769.42999999999
$ws.Range("N138").Value = -110049.43

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2992.5217
$ws.Range("I31").Value = 1664.8
$ws.Range("J31").Value = 4013.8462
$ws.Range("K31").Value = 1664.8
$ws.Range("L31").Value = 4013.8462
$ws.Range("M31").Value = -1369.8
$ws.Range("N31").Value = -4603.8462
$ws.Range("H34").Value = 2992.5217
$ws.Range("I34").Value = 1664.8
$ws.Range("J34").Value = 4013.8462
$ws.Range("K34").Value = 1664.8
$ws.Range("L34").Value = 4013.8462
$ws.Range("M34").Value = -1462.8
$ws.Range("N34").Value = -4417.8462
$ws.Range("H132").Value = 1943
$ws.Range("I132").Value = 1611.8462
$ws.Range("K132").Value = 4835.5386
$ws.Range("M132").Value = -2305.5386
$ws.Range("H134").Value = 2615.6365
$ws.Range("I134").Value = 2494.2942
$ws.Range("J134").Value = 3028.2
$ws.Range("K134").Value = 7482.882599999999
$ws.Range("L134").Value = 9084.599999999999
$ws.Range("M134").Value = -4947.882599999999
$ws.Range("N134").Value = -14154.6
$ws.Range("H138").Value = 94496
$ws.Range("J138").Value = 94496
$ws.Range("L138").Value = 94496
$ws.Range("N138").Value = -104776

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 18058.416
$ws.Range("J93").Value = 18058.416
$ws.Range("L93").Value = 18058.416
$ws.Range("N93").Value = -21802.416
$ws.Range("H122").Value = 254294.67
$ws.Range("I122").Value = 452379
$ws.Range("K122").Value = 1357137
$ws.Range("M122").Value = -1354687
$ws.Range("H126").Value = 3302.0344
$ws.Range("I126").Value = 1577.8462
$ws.Range("K126").Value = 4733.5386
$ws.Range("M126").Value = -2263.5386
$ws.Range("H132").Value = 3212.861
$ws.Range("I132").Value = 2570.913
$ws.Range("J132").Value = 4348.615
$ws.Range("K132").Value = 7712.739
$ws.Range("L132").Value = 13045.845
$ws.Range("M132").Value = -5182.739
$ws.Range("N132").Value = -18105.845
$ws.Range("H135").Value = 56349.816
$ws.Range("J135").Value = 56349.816
$ws.Range("L135").Value = 56349.816
$ws.Range("N135").Value = -66489.81599999999
$ws.Range("H140").Value = 90411.42999999999
$ws.Range("J140").Value = 90396.664
$ws.Range("L140").Value = 90396.664
$ws.Range("N140").Value = -100756.664

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 352134.84
$ws.Range("I68").Value = 352134.84
$ws.Range("K68").Value = 352134.84
$ws.Range("M68").Value = -351385.84
$ws.Range("H71").Value = 352134.84
$ws.Range("I71").Value = 352134.84
$ws.Range("K71").Value = 1760674.2
$ws.Range("M71").Value = -1756930.2
$ws.Range("H122").Value = 85718616
$ws.Range("I122").Value = 142861070
$ws.Range("K122").Value = 428583210
$ws.Range("M122").Value = -428580760
$ws.Range("H136").Value = 5276.1
$ws.Range("I136").Value = 7461
$ws.Range("K136").Value = 22383
$ws.Range("M136").Value = -19833
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2864.8
$ws.Range("I122").Value = 2649.8
$ws.Range("K122").Value = 7949.400000000001
$ws.Range("M122").Value = -5499.400000000001
$ws.Range("H132").Value = 1164.8077
$ws.Range("I132").Value = 999.2857
$ws.Range("K132").Value = 2997.8571
$ws.Range("M132").Value = -467.8571000000002
$ws.Range("H136").Value = 936
$ws.Range("I136").Value = 550.3333
$ws.Range("K136").Value = 1650.9999
$ws.Range("M136").Value = 899.0001
